$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("desc_setor_governo")

$rows = @(
    @(4731, "FUNDO DE DESENVOLVIMENTO DO MINISTÉRIO PÚBLICO", "MINISTÉRIO PÚBLICO", 1090, "PROCURADORIA GERAL DE JUSTIÇA"),
    @(4751, "FUNDO ESPECIAL DA ADVOCACIA GERAL DO ESTADO ", "ADVOCACIA GERAL", 1080, "ADVOCACIA GERAL DO ESTADO"),
    @(4741, "FUNDO ESPECIAL DE GARANTIA DE ACESSO À JUSTIÇA ", "DEFENSORIA PÚBLICA", 1440, "DEFENSORIA PUBLICA DO ESTADO DE MINAS GERAIS"),
    @(2471, "AGÊNCIA REGULADORA DE TRANSPORTES DO ESTADO DE MINAS GERAIS ", "INFRAESTRUTURA, MOBILIDADE E PARCERIAS", 1300, "SECRETARIA DE ESTADO DE INFRAESTRUTURA, MOBILIDADE E PARCERIAS")
)

$startRow = 117
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Cells.Item($r, 1).Value = $data[0]
    $ws.Cells.Item($r, 2).Value = $data[1]
    $ws.Cells.Item($r, 3).Value = $data[2]
    $ws.Cells.Item($r, 4).Value = $data[3]
    $ws.Cells.Item($r, 5).Value = $data[4]
}
